$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, week-covering dates) ---
$ws.Range("A8").Value = "Volume 30   Number  19"
$ws.Range("C9").Value = "Report Covering the Week  5/8/2023  Through  5/14/2023"

# --- Cells changing type (numeric <-> text): copy donor format, then set value ---
$ws.Range("D15").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C15").Value = "0"

$ws.Range("D15").Copy()
$ws.Range("G15").PasteSpecial(-4122)
$ws.Range("G15").Value = "0"

$ws.Range("E15").Copy()
$ws.Range("H15").PasteSpecial(-4122)
$ws.Range("H15").Value = "***.*"

$ws.Range("G17").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("C17").Value = 5

$ws.Range("A18").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("D18").Value = "0"

$ws.Range("A18").Copy()
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("E18").Value = "***.*"

$ws.Range("A20").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("C20").Value = "0"

$ws.Range("A22").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("D22").Value = "0"

$ws.Range("A22").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("E22").Value = "***.*"

$ws.Range("A26").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("C26").Value = "0"

$ws.Range("A26").Copy()
$ws.Range("G26").PasteSpecial(-4122)
$ws.Range("G26").Value = "0"

$ws.Range("A26").Copy()
$ws.Range("H26").PasteSpecial(-4122)
$ws.Range("H26").Value = "***.*"

$excel.CutCopyMode = 0

# --- Simple value updates (same style/type) ---
$ws.Range("N14").Value = 0
$ws.Range("M15").Value = 20
$ws.Range("N15").Value = -45.454545454545
$ws.Range("C16").Value = 2
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 6
$ws.Range("H16").Value = 50
$ws.Range("I16").Value = 30
$ws.Range("J16").Value = 18
$ws.Range("K16").Value = 66.666666666666
$ws.Range("L16").Value = 66.666666666666
$ws.Range("M16").Value = -16.666666666666
$ws.Range("N16").Value = -86.425339366515
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 25
$ws.Range("F17").Value = 10
$ws.Range("H17").Value = 25
$ws.Range("I17").Value = 55
$ws.Range("J17").Value = 40
$ws.Range("K17").Value = 37.5
$ws.Range("L17").Value = 83.333333333333
$ws.Range("M17").Value = 48.648648648648
$ws.Range("N17").Value = -39.560439560439
$ws.Range("F18").Value = 5
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = -28.571428571428
$ws.Range("I18").Value = 45
$ws.Range("K18").Value = 60.714285714285
$ws.Range("L18").Value = 36.363636363636
$ws.Range("M18").Value = -49.438202247191
$ws.Range("N18").Value = -89.510489510489
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = 12.5
$ws.Range("F19").Value = 38
$ws.Range("H19").Value = 31.034482758620
$ws.Range("I19").Value = 154
$ws.Range("J19").Value = 193
$ws.Range("K19").Value = -20.207253886010
$ws.Range("L19").Value = 31.623931623931
$ws.Range("M19").Value = 40
$ws.Range("N19").Value = -14.444444444444
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 10
$ws.Range("G20").Value = 8
$ws.Range("H20").Value = 25
$ws.Range("J20").Value = 37
$ws.Range("K20").Value = 10.810810810810
$ws.Range("L20").Value = 105
$ws.Range("N20").Value = -94.489247311828
$ws.Range("C21").Value = 17
$ws.Range("D21").Value = 14
$ws.Range("E21").Value = 21.428571428571
$ws.Range("F21").Value = 70
$ws.Range("G21").Value = 56
$ws.Range("H21").Value = 25
$ws.Range("I21").Value = 333
$ws.Range("J21").Value = 322
$ws.Range("K21").Value = 3.416149068322
$ws.Range("L21").Value = 50
$ws.Range("M21").Value = -1.479289940828
$ws.Range("N21").Value = -80.154946364719
$ws.Range("L22").Value = 0
$ws.Range("C24").Value = 21
$ws.Range("D24").Value = 36
$ws.Range("E24").Value = -41.666666666666
$ws.Range("F24").Value = 114
$ws.Range("G24").Value = 132
$ws.Range("H24").Value = -13.636363636363
$ws.Range("I24").Value = 580
$ws.Range("J24").Value = 682
$ws.Range("K24").Value = -14.956011730205
$ws.Range("L24").Value = 95.945945945946
$ws.Range("M24").Value = 34.570765661252
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = 28.571428571428
$ws.Range("F25").Value = 28
$ws.Range("G25").Value = 21
$ws.Range("H25").Value = 33.333333333333
$ws.Range("I25").Value = 138
$ws.Range("J25").Value = 123
$ws.Range("K25").Value = 12.195121951219
$ws.Range("L25").Value = 66.265060240963
$ws.Range("M25").Value = 7.8125
$ws.Range("C27").Value = 3
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 200
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = 8
$ws.Range("H27").Value = -50
$ws.Range("I27").Value = 9
$ws.Range("J27").Value = 17
$ws.Range("K27").Value = -47.058823529411
$ws.Range("L27").Value = 28.571428571428
